$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.108.97'
$ws.Range("E2").Value = '  +6.49%  '
$ws.Range("D3").Value = '2.638.24'
$ws.Range("E3").Value = '  +10.43%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.996'
$cell.ClearFormats()
$ws.Range("E4").Value = '  -0.49%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '313.54'
$cell.ClearFormats()
$ws.Range("E5").Value = '  +7.12%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '105.98'
$cell.ClearFormats()
$ws.Range("E6").Value = '  +13.25%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.613'
$cell.ClearFormats()
$ws.Range("E7").Value = '  +10.34%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.ClearFormats()
$ws.Range("E8").Value = '  -0.22%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.599'
$cell.ClearFormats()
$ws.Range("E9").Value = '  +20.06%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '40.26'
$cell.ClearFormats()
$ws.Range("E10").Value = '  +17.84%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0865'
$cell.ClearFormats()
$ws.Range("E11").Value = '  +11.55%  '
$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '55.41'
$cell.ClearFormats()
$ws.Range("E12").Value = '  +3.79%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '8.41'
$cell.ClearFormats()
$ws.Range("E13").Value = '  +20.71%  '
$ws.Range("D14").Value = '3.030.50'
$ws.Range("E14").Value = '  +10.07%  '
$ws.Range("E15").Value = '  +3.34%  '
$ws.Range("D16").Value = '2.640.42'
$ws.Range("E16").Value = '  +10.98%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.941'
$cell.ClearFormats()
$ws.Range("E17").Value = '  +14.23%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '15.35'
$cell.ClearFormats()
$ws.Range("E18").Value = '  +8.94%  '
$ws.Range("D19").Value = '47.970.25'
$ws.Range("E19").Value = '  +6.21%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.0000104'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +11.47%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '13.42'
$cell.ClearFormats()
$ws.Range("E21").Value = '  +8.24%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.86'
$cell.ClearFormats()
$ws.Range("E22").Value = '  +13.01%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '73.28'
$cell.ClearFormats()
$ws.Range("E23").Value = '  +10.33%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '276.59'
$cell.ClearFormats()
$ws.Range("E24").Value = '  +16.35%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '3.11'
$cell.ClearFormats()
$ws.Range("E25").Value = '  +13.03%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.24'
$cell.ClearFormats()
$ws.Range("E26").Value = '  +19.25%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '30.64'
$cell.ClearFormats()
$ws.Range("E27").Value = '  +46.75%  '
$ws.Range("E28").Value = '  +0.13%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '4.11'
$cell.ClearFormats()
$ws.Range("E29").Value = '  +2.17%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '10.76'
$cell.ClearFormats()
$ws.Range("E30").Value = '  +13.05%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '2.32'
$cell.ClearFormats()
$ws.Range("E31").Value = '  +4.75%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '40.28'
$cell.ClearFormats()
$ws.Range("E32").Value = '  +8.13%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '6.23'
$cell.ClearFormats()
$ws.Range("E33").Value = '  +15.70%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '3.73'
$cell.ClearFormats()
$ws.Range("E34").Value = '  -1.83%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.0865'
$cell.ClearFormats()
$ws.Range("E35").Value = '  +14.40%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '2.89'
$cell.ClearFormats()
$ws.Range("E36").Value = '  +7.02%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.23'
$cell.ClearFormats()
$ws.Range("E37").Value = '  +13.95%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '152.48'
$cell.ClearFormats()
$ws.Range("E38").Value = '  +3.85%  '
$ws.Range("E39").Value = '  +12.26%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.125'
$cell.ClearFormats()
$ws.Range("E40").Value = '  +9.87%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '16.62'
$cell.ClearFormats()
$ws.Range("E41").Value = '  +15.16%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '4.35'
$cell.ClearFormats()
$ws.Range("E42").Value = '  +18.47%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '22.77'
$cell.ClearFormats()
$ws.Range("E43").Value = '  +49.42%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '3.73'
$cell.ClearFormats()
$ws.Range("E44").Value = '  +18.82%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.0336'
$cell.ClearFormats()
$ws.Range("E45").Value = '  +15.21%  '
$ws.Range("D46").Value = '2.204.84'
$ws.Range("E46").Value = '  +12.04%  '
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '97.77'
$cell.ClearFormats()
$ws.Range("E47").Value = '  +11.08%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$ws.Range("E48").Value = '  -0.01%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '9.99'
$cell.ClearFormats()
$ws.Range("E49").Value = '  +19.26%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '115.06'
$cell.ClearFormats()
$ws.Range("E50").Value = '  +16.18%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.87'
$cell.ClearFormats()
$ws.Range("E51").Value = '  +10.34%  '
